$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Data"
$ws.Range("A1").Value = "`"- Kitchen Assistant"
$ws.Range("A2").Value = "Dishwasher/Prep Cook`""
